$d = $word.ActiveDocument

# The new paragraphs are plain, empty paragraphs whose only formatting is
# the Portuguese (Brazil) language mark carried on the paragraph mark's
# run properties - matching the other blank paragraphs already present
# in this document.
$paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p>'

# Append two such empty paragraphs at the very end of the document body
# (right before the sectPr), after the existing "nc -u -l 6664" paragraph.
for ($i = 0; $i -lt 2; $i++) {
    $endPos = $d.Content.End
    $r = $d.Range($endPos, $endPos)
    [void]$r.InsertXML($paraXml)
}
